# Regenerate the quadratic/linear problem data (new random experiment values).
# NOTE: worksheets "Vector_bf" and "Vector_BF" differ only by letter case,
# so we must address sheets by their 1-based index (Worksheets.Item($n))
# rather than by name, to avoid ambiguous case-insensitive name lookups.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$val) {
    # Force the cell to be stored as text (string) even though the value
    # looks like a plain number, then strip the temporary style back off
    # so no extraneous formatting is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# --- Sheet 3: Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Range("A2") "5.185921956581478 - x + 0.29310250068700194y_1 - 0.11426215993404765y_2"
Set-TextValue $ws3.Range("B2") "-5.185921956581478"
Set-TextValue $ws3.Range("D2") "0.34"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "0.2"

Set-TextValue $ws3.Range("A3") "-5.591412109628146 - 0.05181107677306451y_1 + 1.4946028789848305y_2"
Set-TextValue $ws3.Range("B3") "5.591412109628146"
Set-TextValue $ws3.Range("D3") "0.14"
Set-TextValue $ws3.Range("E3") "-8.5"
Set-TextValue $ws3.Range("F3") "-1.0"

Set-TextValue $ws3.Range("A4") "-2.7820335256938713 + 0.7766144545204726y_1 + 0.22415498763396524y_2"
Set-TextValue $ws3.Range("B4") "2.1720335256938714"
Set-TextValue $ws3.Range("D4") "0.38"
Set-TextValue $ws3.Range("E4") "5.4"
Set-TextValue $ws3.Range("F4") "2.0"

# --- Sheet 4: Punto_modificado ---
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4.Range("A2") "5.25"
Set-TextValue $ws4.Range("B2") "1.7000000000000002"
Set-TextValue $ws4.Range("C2") "3.8"

# --- Sheet 5: Vector_bf ---
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5.Range("A2") "-2.087514792203131"
Set-TextValue $ws5.Range("A3") "-1.2555741639812068"

# --- Sheet 6: Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6.Range("A2") "-0.0"
Set-TextValue $ws6.Range("A3") "-6.1041122069816005"
Set-TextValue $ws6.Range("A4") "9.693687538147646"

# --- Sheet 7: Vector_Alpha ---
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("A2").Value = 1.77
$ws7.Range("A3").Value = 2.79
